$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1772.6364
$ws.Range("J17").Value = 1772.6364
$ws.Range("L17").Value = 5317.9092
$ws.Range("N17").Value = -5653.9092

$ws.Range("H74").Value = 3666.6667

$ws.Range("H77").Value = 3666.6667

$ws.Range("H98").Value = 55556572
$ws.Range("I98").Value = 55556572
$ws.Range("K98").Value = 55556572
$ws.Range("M98").Value = -55555074

$ws.Range("H113").Value = 125001380
$ws.Range("I113").Value = 50001250
$ws.Range("K113").Value = 50001250
$ws.Range("M113").Value = -49997996

$ws.Range("H122").Value = 55556572
$ws.Range("I122").Value = 55556572
$ws.Range("K122").Value = 166669716
$ws.Range("M122").Value = -166667266

$ws.Range("H129").Value = 1463.1666
$ws.Range("I129").Value = 895.55554
$ws.Range("J129").Value = 3166
$ws.Range("K129").Value = 2686.66662
$ws.Range("L129").Value = 9498
$ws.Range("M129").Value = 2313.33338
$ws.Range("N129").Value = -19498

$ws.Range("H132").Value = 945.78125
$ws.Range("I132").Value = 831.129
$ws.Range("K132").Value = 2493.387
$ws.Range("M132").Value = 36.61299999999983

$ws.Range("H135").Value = 3503.6
$ws.Range("I135").Value = 2827.3333
$ws.Range("K135").Value = 25445.9997
$ws.Range("M135").Value = -22910.9997

$ws.Range("H137").Value = 3489
$ws.Range("I137").Value = 2307.5227
$ws.Range("K137").Value = 6922.5681
$ws.Range("M137").Value = -4372.5681

$ws.Range("H138").Value = 2128.6104
$ws.Range("I138").Value = 722.4138
$ws.Range("J138").Value = 2978.1875
$ws.Range("K138").Value = 2167.2414
$ws.Range("L138").Value = 8934.5625
$ws.Range("M138").Value = 2972.7586
$ws.Range("N138").Value = -19214.5625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11910798
$ws.Range("I32").Value = 16131732
$ws.Range("K32").Value = 16131732
$ws.Range("M32").Value = -16131445

$ws.Range("H45").Value = 33336254
$ws.Range("I45").Value = 62501830
$ws.Range("J45").Value = 4167
$ws.Range("K45").Value = 62501830
$ws.Range("L45").Value = 4167
$ws.Range("M45").Value = -62501453
$ws.Range("N45").Value = -4921

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H2").Value = 81359.8
$ws.Range("J2").Value = 81359.8
$ws.Range("L2").Value = 81359.8
$ws.Range("N2").Value = -81585.8

$ws.Range("H11").Value = 1000.4
$ws.Range("I11").Value = 1002
$ws.Range("J11").Value = 999.3333
$ws.Range("K11").Value = 1002
$ws.Range("L11").Value = 999.3333
$ws.Range("M11").Value = -862
$ws.Range("N11").Value = -1279.3333

$ws.Range("H86").Value = 2642.8635
$ws.Range("J86").Value = 1823.5714
$ws.Range("L86").Value = 1823.5714
$ws.Range("N86").Value = -4069.5714

$ws.Range("H89").Value = 2642.8635
$ws.Range("J89").Value = 1823.5714
$ws.Range("L89").Value = 9117.857
$ws.Range("N89").Value = -20349.857

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 238.25
$ws.Range("J7").Value = 223.25
$ws.Range("L7").Value = 223.25
$ws.Range("N7").Value = -449.25

$ws.Range("H31").Value = 711069.8
$ws.Range("I31").Value = 3140.3125
$ws.Range("J31").Value = 1377356.4
$ws.Range("K31").Value = 3140.3125
$ws.Range("L31").Value = 1377356.4
$ws.Range("M31").Value = -2845.3125
$ws.Range("N31").Value = -1377946.4

$ws.Range("H34").Value = 711069.8
$ws.Range("I34").Value = 3140.3125
$ws.Range("J34").Value = 1377356.4
$ws.Range("K34").Value = 3140.3125
$ws.Range("L34").Value = 1377356.4
$ws.Range("M34").Value = -2938.3125
$ws.Range("N34").Value = -1377760.4

$ws.Range("H58").Value = 2158.16
$ws.Range("J58").Value = 6951
$ws.Range("L58").Value = 6951
$ws.Range("N58").Value = -7357

$ws.Range("H86").Value = 8459.25
$ws.Range("I86").Value = 9224.666999999999
$ws.Range("J86").Value = 8000
$ws.Range("K86").Value = 9224.666999999999
$ws.Range("L86").Value = 8000
$ws.Range("M86").Value = -8101.666999999999
$ws.Range("N86").Value = -10246

$ws.Range("H89").Value = 8459.25
$ws.Range("I89").Value = 9224.666999999999
$ws.Range("J89").Value = 8000
$ws.Range("K89").Value = 46123.335
$ws.Range("L89").Value = 40000
$ws.Range("M89").Value = -40507.335
$ws.Range("N89").Value = -51232

$ws.Range("H98").Value = 42192.25
$ws.Range("J98").Value = 54256.332
$ws.Range("L98").Value = 54256.332
$ws.Range("N98").Value = -58748.332

$ws.Range("H99").Value = 3591.9092
$ws.Range("I99").Value = 3296.1428
$ws.Range("J99").Value = 4109.5
$ws.Range("K99").Value = 3296.1428
$ws.Range("L99").Value = 4109.5
$ws.Range("M99").Value = -1798.1428
$ws.Range("N99").Value = -7105.5

$ws.Range("H107").Value = 2016.6786
$ws.Range("I107").Value = 813.3333
$ws.Range("K107").Value = 813.3333
$ws.Range("M107").Value = 1106.6667

$ws.Range("H122").Value = 2136.5
$ws.Range("I122").Value = 1821.7693
$ws.Range("K122").Value = 5465.3079
$ws.Range("M122").Value = -3015.3079

$ws.Range("H126").Value = 3591.9092
$ws.Range("I126").Value = 3296.1428
$ws.Range("J126").Value = 4109.5
$ws.Range("K126").Value = 9888.428400000001
$ws.Range("L126").Value = 12328.5
$ws.Range("M126").Value = -7418.428400000001
$ws.Range("N126").Value = -17268.5

$ws.Range("H136").Value = 2158.16
$ws.Range("J136").Value = 6951
$ws.Range("L136").Value = 20853
$ws.Range("N136").Value = -25953

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 589.6896400000001
$ws.Range("J107").Value = 732.3333
$ws.Range("L107").Value = 2196.9999
$ws.Range("N107").Value = -6036.9999

$ws.Range("H131").Value = 11202.451
$ws.Range("J131").Value = 11202.451
$ws.Range("L131").Value = 33607.353
$ws.Range("N131").Value = -43687.353

$ws.Range("H132").Value = 2124.5
$ws.Range("I132").Value = 1000
$ws.Range("J132").Value = 2499.3333
$ws.Range("K132").Value = 9000
$ws.Range("L132").Value = 22493.9997
$ws.Range("M132").Value = -6470
$ws.Range("N132").Value = -27553.9997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 28000
$ws.Range("J15").Value = 28000
$ws.Range("L15").Value = 28000
$ws.Range("N15").Value = -28576

$ws.Range("H80").Value = 4939.4614
$ws.Range("I80").Value = 3629.6
$ws.Range("J80").Value = 5758.125
$ws.Range("K80").Value = 3629.6
$ws.Range("L80").Value = 5758.125
$ws.Range("M80").Value = -2631.6
$ws.Range("N80").Value = -7754.125

$ws.Range("H81").Value = 28000
$ws.Range("J81").Value = 28000
$ws.Range("L81").Value = 28000
$ws.Range("N81").Value = -29996

$ws.Range("H83").Value = 4939.4614
$ws.Range("I83").Value = 3629.6
$ws.Range("J83").Value = 5758.125
$ws.Range("K83").Value = 18148
$ws.Range("L83").Value = 28790.625
$ws.Range("M83").Value = -13156
$ws.Range("N83").Value = -38774.625

$ws.Range("H84").Value = 28000
$ws.Range("J84").Value = 28000
$ws.Range("L84").Value = 84000
$ws.Range("N84").Value = -93984

$ws.Range("I122").Value = 1444.9333
$ws.Range("J122").Value = 1177
$ws.Range("K122").Value = 4334.7999
$ws.Range("L122").Value = 3531
$ws.Range("M122").Value = -1884.7999
$ws.Range("N122").Value = -8431

$ws.Range("H126").Value = 4438.385
$ws.Range("I126").Value = 3927.75
$ws.Range("J126").Value = 4665.3335
$ws.Range("K126").Value = 11783.25
$ws.Range("L126").Value = 13996.0005
$ws.Range("M126").Value = -9313.25
$ws.Range("N126").Value = -18936.0005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1599.3636
$ws.Range("I16").Value = 1359.3
$ws.Range("K16").Value = 1359.3
$ws.Range("M16").Value = -1189.3

$ws.Range("H40").Value = 5363.636
$ws.Range("I40").Value = 5000
$ws.Range("J40").Value = 5444.4443
$ws.Range("K40").Value = 5000
$ws.Range("L40").Value = 5444.4443
$ws.Range("M40").Value = -4864
$ws.Range("N40").Value = -5716.4443

$ws.Range("H93").Value = 500001000
$ws.Range("I93").Value = 1000000000
$ws.Range("J93").Value = 2004
$ws.Range("K93").Value = 1000000000
$ws.Range("L93").Value = 2004
$ws.Range("M93").Value = -999998752
$ws.Range("N93").Value = -4500

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 7260
$ws.Range("I81").Value = 2283.3333
$ws.Range("J81").Value = 14725
$ws.Range("K81").Value = 4566.6666
$ws.Range("L81").Value = 29450
$ws.Range("M81").Value = -3505.6666
$ws.Range("N81").Value = -31572

$ws.Range("H84").Value = 7260
$ws.Range("I84").Value = 2283.3333
$ws.Range("J84").Value = 14725
$ws.Range("K84").Value = 22833.333
$ws.Range("L84").Value = 147250
$ws.Range("M84").Value = -17529.333
$ws.Range("N84").Value = -157858

$ws.Range("H107").Value = 13158811
$ws.Range("I107").Value = 17858172
$ws.Range("J107").Value = 599.2
$ws.Range("K107").Value = 53574516
$ws.Range("L107").Value = 1797.6
$ws.Range("M107").Value = -53572596
$ws.Range("N107").Value = -5637.6

$ws.Range("H122").Value = 2370.875
$ws.Range("I122").Value = 2273.4736
$ws.Range("K122").Value = 6820.4208
$ws.Range("M122").Value = -4370.4208

$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()
